# Daily attendance processing - 2025-10-26 05:20:12
# Reverses the order of names/emails in the "Recorded By" column (G) for the
# specific attendance rows that were re-processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G ("Recorded By") whose comma-separated list of recorders
# needs to be reversed in order.
$targetRows = @(2, 4, 5, 7, 8, 11, 17, 29, 31, 32, 34, 35, 38, 44, 56, 58, 59, 61, 62, 65, 71, 83, 84, 85, 90, 96, 97, 99, 109, 110, 111, 116, 122, 123, 125, 135, 136, 137, 142, 148, 149, 151)

foreach ($row in $targetRows) {
    $cell = $ws.Range("G$row")
    $current = [string]$cell.Value2

    # Split on comma, trim whitespace from each part, then reverse the order.
    $parts = $current.Split(",") | ForEach-Object { $_.Trim() }
    $reversed = $parts[($parts.Count - 1)..0]

    $cell.Value = [string]::Join(", ", $reversed)
}
